{"js": "// The document originally starts with two paragraphs that each contain a\n// single inline picture, followed by a paragraph of syntax-highlighted C#\n// code. The code paragraph holds a collapsed \"_GoBack\" bookmark placed\n// right after the opening \"{\" of Main() (before the Console.WriteLine\n// call). The edit:\n//   1) removes the two leading picture paragraphs entirely, and\n//   2) relocates the \"_GoBack\" bookmark to the very start of the code\n//      paragraph (now the first paragraph), before the \"namespace\" text.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// 1) Delete the two leading paragraphs that only contain a picture.\nconst picParagraphs = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  p.inlinePictures.load(\"items\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.inlinePictures.items.length > 0) {\n    picParagraphs.push(p);\n  }\n}\npicParagraphs.forEach((p) => p.delete());\nawait context.sync();\n\n// 2) Move the \"_GoBack\" bookmark to the start of the (now first)\n//    paragraph, ahead of the \"namespace\" run.\nconst refreshedParagraphs = body.paragraphs;\nrefreshedParagraphs.load(\"items\");\nawait context.sync();\n\nconst targetParagraph = refreshedParagraphs.items[0];\nconst startRange = targetParagraph.getRange(\"Start\");\n\ncontext.document.deleteBookmark(\"_GoBack\");\nstartRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The document starts with two paragraphs that each hold a single inline\n# picture, followed by a paragraph of syntax-highlighted C# code. That code\n# paragraph carries a collapsed \"_GoBack\" bookmark placed right after the\n# opening \"{\" of Main() (before the Console.WriteLine call). This script:\n#   1) relocates the \"_GoBack\" bookmark to the very start of the code\n#      paragraph (before the \"namespace\" text) while the picture paragraphs\n#      are still present (so the target position is not a paragraph-start\n#      boundary at insertion time, keeping the bookmark collapsed), then\n#   2) removes the two leading picture-only paragraphs entirely.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that contains the C# code (first paragraph without\n# an inline picture).\n$codeParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.InlineShapes.Count -eq 0) {\n        $codeParagraph = $p\n        break\n    }\n}\n\n# Move the \"_GoBack\" bookmark to the start of that paragraph.\n$targetStart = $codeParagraph.Range.Start\n$targetRange = $d.Range($targetStart, $targetStart)\n\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n$d.Bookmarks.Add(\"_GoBack\", $targetRange)\n\n# Delete the leading paragraphs that only contain a picture (walk back to\n# front so the indices of paragraphs still to be removed stay valid).\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.InlineShapes.Count -gt 0) {\n        $p.Range.Delete()\n    }\n}\n"}
